# C5-PowerPoint.pptx edit:
#  1) The "Sources of finance" table on slide 6 is switched to a
#     different built-in table style.
#  2) The deck's theme colour scheme is switched from the "Integral"
#     palette to the standard "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{E8E6D363-C262-48A2-BA01-684D8A2BB90C}")
    }
}

# --- 2. Theme colours -------------------------------------------------
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

function Set-ThemeRgb {
    param($Scheme, $Index, $Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $Scheme.Item($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office Theme colours, in MsoThemeColorSchemeIndex order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
Set-ThemeRgb $colorScheme 1  "000000"
Set-ThemeRgb $colorScheme 2  "FFFFFF"
Set-ThemeRgb $colorScheme 3  "44546A"
Set-ThemeRgb $colorScheme 4  "E7E6E6"
Set-ThemeRgb $colorScheme 5  "5B9BD5"
Set-ThemeRgb $colorScheme 6  "ED7D31"
Set-ThemeRgb $colorScheme 7  "A5A5A5"
Set-ThemeRgb $colorScheme 8  "FFC000"
Set-ThemeRgb $colorScheme 9  "4472C4"
Set-ThemeRgb $colorScheme 10 "70AD47"
Set-ThemeRgb $colorScheme 11 "0563C1"
Set-ThemeRgb $colorScheme 12 "954F72"
